$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.9
$ws.Range("A3").Value = -21.557
$ws.Range("E5").Value = 13.098
$ws.Range("A14").Value = -21.04
$ws.Range("A16").Value = -21.008
$ws.Range("E16").Value = 12.963
$ws.Range("A21").Value = -21.04
$ws.Range("A23").Value = -21.584
$ws.Range("A25").Value = -22.27
